# End of Year Wifcaster Update
# - Rename a handful of NPC "Events" entries.
# - Add a new header/number row (row 1) with sequence numbers 1-9 across B1:J1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text renames (Events column, D, except the one NPC name in C6) ---
$ws.Range("D3").Value  = "Splash"          # was "River"
$ws.Range("D4").Value  = "Brine"           # was "Basin"
$ws.Range("C6").Value  = "Cyran Herqen"    # was "Ralnor Beifiel"
$ws.Range("D8").Value  = "Quartz"          # was "Slab"
$ws.Range("D14").Value = "Kindle"          # was "Cinder"
$ws.Range("D18").Value = "Flurry"          # was "Drift"

# --- New row 1: sequence numbers 1-9 in columns B..J, formatted like the ---
# --- existing "no fill" body cells (e.g. B3) rather than the shaded      ---
# --- header row (row 2) or the tinted Events column (D).                ---
$ws.Range("B3").Copy()
$ws.Range("B1:J1").PasteSpecial(-4122)

$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
